# Fruta / hortaliza, semanal
# Insert 2 new rows before row 226 (pushing old rows 226-263 down to 228-265)
# and populate them with the new week's data for
# "Vega Monumental Concepción - Limón".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("226:227").Insert()

# Row 226 - 1a amarillo
$ws.Cells.Item(226, 1).Value = 11
$ws.Cells.Item(226, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(226, 3).Value = "Bíobío"
$ws.Cells.Item(226, 4).Value = 44474
$ws.Cells.Item(226, 5).Value = 8
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100102
$ws.Cells.Item(226, 8).Value = "Cítricos"
$ws.Cells.Item(226, 9).Value = 100102003
$ws.Cells.Item(226, 10).Value = "Limón"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "1a amarillo"
$ws.Cells.Item(226, 13).Value = 400
$ws.Cells.Item(226, 14).Value = 6500
$ws.Cells.Item(226, 15).Value = 6500
$ws.Cells.Item(226, 16).Value = 6500
$ws.Cells.Item(226, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(226, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(226, 19).Value = 406
$ws.Cells.Item(226, 20).Value = 16

# Row 227 - 2a amarillo
$ws.Cells.Item(227, 1).Value = 11
$ws.Cells.Item(227, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(227, 3).Value = "Bíobío"
$ws.Cells.Item(227, 4).Value = 44474
$ws.Cells.Item(227, 5).Value = 8
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100102
$ws.Cells.Item(227, 8).Value = "Cítricos"
$ws.Cells.Item(227, 9).Value = 100102003
$ws.Cells.Item(227, 10).Value = "Limón"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "2a amarillo"
$ws.Cells.Item(227, 13).Value = 400
$ws.Cells.Item(227, 14).Value = 5500
$ws.Cells.Item(227, 15).Value = 5500
$ws.Cells.Item(227, 16).Value = 5500
$ws.Cells.Item(227, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(227, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(227, 19).Value = 344
$ws.Cells.Item(227, 20).Value = 16

# Ensure column D keeps the date style used throughout the column.
$ws.Range("D226:D227").NumberFormat = $ws.Range("D228").NumberFormat
